# Regenerate orders with updated distance/sizes.
# Applies text substitutions across the used range:
#   D51 -> D55
#   D64 -> D69
#   D80 -> D86
#   S30 -> S31
# These substitutions affect the Condition, Filename_Left, Filename_Right,
# Distance and Size columns (and any string that embeds those tokens).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count
$startRow = $used.Row
$startCol = $used.Column

for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $cell = $ws.Cells.Item($startRow + $r, $startCol + $c)
        $val = $cell.Value2

        if ($val -is [string]) {
            $newVal = $val
            $newVal = $newVal.Replace("D51", "D55")
            $newVal = $newVal.Replace("D64", "D69")
            $newVal = $newVal.Replace("D80", "D86")
            $newVal = $newVal.Replace("S30", "S31")

            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
